# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" sheets to the newly scraped numbers.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet (row => new F value) ---
$exhibitUpdates = @{
    2  = 2722
    5  = 1520
    6  = 1143
    8  = 540
    12 = 9225
    18 = 471
    19 = 636
    21 = 1182
    23 = 2097
    24 = 2197
    26 = 1902
    28 = 1931
    31 = 282
    33 = 213
    34 = 27
    37 = 299
    40 = 64
    41 = 779
    42 = 49
    43 = 1401
    44 = 305
    46 = 183
    47 = 654
    48 = 79
    49 = 300
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# --- 全部类型 sheet (row => new F value) ---
$allUpdates = @{
    2  = 2722
    4  = 1520
    6  = 1143
    8  = 540
    10 = 9225
    18 = 471
    19 = 636
    20 = 1182
    22 = 2197
    23 = 1902
    26 = 282
    28 = 213
    29 = 27
    32 = 299
    38 = 64
    39 = 779
    41 = 49
    42 = 1401
    44 = 305
    46 = 183
    47 = 654
    48 = 300
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
